$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.204.87'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").Value = '3.424.72'
$ws.Range("E3").Value = '  +1.58%  '

$ws.Range("E4").Value = '  -0.02%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '572.90'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.21%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '139.19'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.84%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '3.426.24'
$ws.Range("E8").Value = '  +1.61%  '

$ws.Range("E9").Value = '  -0.31%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '7.69'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.86%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.122'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.35%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.383'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.88%  '

$ws.Range("D13").Value = '4.004.88'
$ws.Range("E13").Value = '  +1.43%  '

$ws.Range("E14").Value = '  -0.83%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '26.69'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +3.03%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.0000173'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.07%  '

$ws.Range("D17").Value = '3.423.04'
$ws.Range("E17").Value = '  +1.50%  '

$ws.Range("D18").Value = '61.316.59'
$ws.Range("E18").Value = '  -0.24%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '5.95'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.40%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.94'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.25%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '9.42'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.01%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '379.79'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +1.07%  '

$ws.Range("B23").Value = 'WrappedeETH'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D23").Value = '3.556.47'
$ws.Range("E23").Value = '  +1.21%  '

$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.555'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.31%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.15%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '71.52'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("E27").Value = '  -2.12%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.177'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +9.77%  '

$ws.Range("E29").Value = '  -7.42%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '7.58'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.10%  '

$ws.Range("E31").Value = '  -0.21%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '8.16'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.69%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '2.15'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.13%  '

$ws.Range("E34").Value = '  -0.02%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '23.74'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.01%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '5.20'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.10%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '6.95'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +2.15%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.57'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +2.42%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '166.45'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.81%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0779'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.08%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '26.52'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +7.74%  '

$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.783'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.40%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.73'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.77%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '41.97'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.19%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '4.41'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.83%  '

$ws.Range("E47").Value = '  -1.97%  '

$ws.Range("D48").Value = '2.607.76'
$ws.Range("E48").Value = '  +10.60%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '23.86'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +5.51%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '6.79'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.48%  '

$ws.Range("E51").Value = '  -0.49%  '
